$d = $word.ActiveDocument

# Change 1 (cell "Afeta by" -> "O problema" description): the long
# original sentence describing the problem is replaced with the much
# shorter "falta de comunicação."
$ok1 = $d.Content.Find.Execute(
    "Do salão é o atendimento e a falta de comunicação entre o quadro de funcionários e perda de clientes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "falta de comunicação.", 2)

# Change 2 ("Afeta" row): "O proprietário, Funcionários e clientes" ->
# "o proprietário, cabelereiros, manicures, esteticista e clientes."
$ok2 = $d.Content.Find.Execute(
    "O proprietário, Funcionários e clientes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "o proprietário, cabelereiros, manicures, esteticista e clientes.", 2)

# Change 3 ("Devido" row): the old text about control of schedules /
# services / billing confusion is replaced with a short note about poor
# communication between staff.
$ok3 = $d.Content.Find.Execute(
    "A falta de controle dos horários de atendimento, serviços e  confusão no momento da cobrança.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ao desempenho comunicação ruim entre os colaboradores ", 2)

# Change 4 ("Os benefícios desse" row): "Sistema " -> "Sistema de
# agendamento  " right before "São:" (keeps the existing color formatting
# on the run(s) it touches).
$ok4 = $d.Content.Find.Execute(
    "Sistema São:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sistema de agendamento  São:", 2)

Write-Host "Find/Replace results: $ok1 $ok2 $ok3 $ok4"
